# 2023_summary_data.xlsx -- "updated for corrected data" edit
#
# statements_combined updated to remove duplicated data points (by summing).
# This adjusts six raw counts on "By New Statement Type" (column C); the
# "By New OLE Category" sheet's SUM()/ratio formulas recalculate on their own
# since they reference those cells. Sheet selections / active-tab bookkeeping
# are updated to match where the author was last working (the
# "By New Statement Type" sheet, near the bottom of the data).

$wb = $excel.ActiveWorkbook

$wsNewType = $wb.Worksheets.Item("By New Statement Type")
$wsNewCat  = $wb.Worksheets.Item("By New OLE Category")

# --- corrected raw counts on "By New Statement Type" ---------------------
$wsNewType.Range("C3").Value  = 8
$wsNewType.Range("C5").Value  = 11
$wsNewType.Range("C6").Value  = 16
$wsNewType.Range("C35").Value = 14
$wsNewType.Range("C39").Value = 24
$wsNewType.Range("C46").Value = 18

# "By New OLE Category" (sheet2) B/C columns are formulas
# (=SUM('By New Statement Type'!...) and =Bn/B14) and recalculate
# automatically from the edits above.

# --- view/selection state --------------------------------------------------
# Update the selection remembered on "By New OLE Category" (was B23).
$wsNewCat.Activate()
$wsNewCat.Range("F16").Select()

# "By New Statement Type" ends up the active sheet/tab (this also clears
# tabSelected from whichever sheet -- "By Old OLE Category (2023)" -- had it
# before), scrolled near the bottom of the data, with C46 selected.
$wsNewType.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$wsNewType.Range("C46").Select()
